$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.320746
$ws.Range("H2").Value = 288.962238
$ws.Range("I2").Value = 0.3809824610908788
$ws.Range("J2").Value = 0.3809824610908788
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 19.816421
$ws.Range("N2").Value = 59.449263
$ws.Range("O2").Value = 0.4265240049466206
$ws.Range("P2").Value = 0.4265240049466206
$ws.Range("Q2").Value = 1908.732453770066
$ws.Range("R2").Value = 17178.59208393059
$ws.Range("S2").Value = 0.1624981651189017
$ws.Range("T2").Value = 0.1624981651189017
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.320746
$ws.Range("H3").Value = 288.962238
$ws.Range("I3").Value = 0.3809824610908788
$ws.Range("J3").Value = 0.3809824610908788
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.17573333333333
$ws.Range("N3").Value = 54.5272
$ws.Range("O3").Value = 0.3912102278294916
$ws.Range("P3").Value = 0.3912102278294917
$ws.Range("Q3").Value = 1750.700193763733
$ws.Range("R3").Value = 15756.3017438736
$ws.Range("S3").Value = 0.1490442354024031
$ws.Range("T3").Value = 0.1490442354024032
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.320746
$ws.Range("H4").Value = 288.962238
$ws.Range("I4").Value = 0.3809824610908788
$ws.Range("J4").Value = 0.3809824610908788
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.429072666666666
$ws.Range("N4").Value = 19.287218
$ws.Range("O4").Value = 0.1383778545015528
$ws.Range("P4").Value = 0.1383778545015528
$ws.Range("Q4").Value = 619.2530753415426
$ws.Range("R4").Value = 5573.277678073884
$ws.Range("S4").Value = 0.05271953556847714
$ws.Range("T4").Value = 0.05271953556847715
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.320746
$ws.Range("H5").Value = 288.962238
$ws.Range("I5").Value = 0.3809824610908788
$ws.Range("J5").Value = 0.3809824610908788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.039044333333333
$ws.Range("N5").Value = 6.117133
$ws.Range("O5").Value = 0.04388791272233494
$ws.Range("P5").Value = 0.04388791272233494
$ws.Range("Q5").Value = 196.4022713137393
$ws.Range("R5").Value = 1767.620441823654
$ws.Range("S5").Value = 0.01672052500109686
$ws.Range("T5").Value = 0.01672052500109686
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.46467533333333
$ws.Range("H6").Value = 55.394026
$ws.Range("I6").Value = 0.07303429161291354
$ws.Range("J6").Value = 0.07303429161291354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.816421
$ws.Range("N6").Value = 59.449263
$ws.Range("O6").Value = 0.4265240049466206
$ws.Range("P6").Value = 0.4265240049466206
$ws.Range("Q6").Value = 365.9037800336487
$ws.Range("R6").Value = 3293.134020302838
$ws.Range("S6").Value = 0.03115087855717926
$ws.Range("T6").Value = 0.03115087855717927
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.46467533333333
$ws.Range("H7").Value = 55.394026
$ws.Range("I7").Value = 0.07303429161291354
$ws.Range("J7").Value = 0.07303429161291354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.17573333333333
$ws.Range("N7").Value = 54.5272
$ws.Range("O7").Value = 0.3912102278294916
$ws.Range("P7").Value = 0.3912102278294917
$ws.Range("Q7").Value = 335.6090149452444
$ws.Range("R7").Value = 3020.4811345072
$ws.Range("S7").Value = 0.02857176186125343
$ws.Range("T7").Value = 0.02857176186125344
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.46467533333333
$ws.Range("H8").Value = 55.394026
$ws.Range("I8").Value = 0.07303429161291354
$ws.Range("J8").Value = 0.07303429161291354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.429072666666666
$ws.Range("N8").Value = 19.287218
$ws.Range("O8").Value = 0.1383778545015528
$ws.Range("P8").Value = 0.1383778545015528
$ws.Range("Q8").Value = 118.7107394844075
$ws.Range("R8").Value = 1068.396655359668
$ws.Range("S8").Value = 0.01010632857843573
$ws.Range("T8").Value = 0.01010632857843573
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.46467533333333
$ws.Range("H9").Value = 55.394026
$ws.Range("I9").Value = 0.07303429161291354
$ws.Range("J9").Value = 0.07303429161291354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.039044333333333
$ws.Range("N9").Value = 6.117133
$ws.Range("O9").Value = 0.04388791272233494
$ws.Range("P9").Value = 0.04388791272233494
$ws.Range("Q9").Value = 37.65029160527311
$ws.Range("R9").Value = 338.852624447458
$ws.Range("S9").Value = 0.003205322616045108
$ws.Range("T9").Value = 0.003205322616045108
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.909391
$ws.Range("H10").Value = 368.728173
$ws.Range("I10").Value = 0.4861499128584522
$ws.Range("J10").Value = 0.4861499128584522
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.816421
$ws.Range("N10").Value = 59.449263
$ws.Range("O10").Value = 0.4265240049466206
$ws.Range("P10").Value = 0.4265240049466206
$ws.Range("Q10").Value = 2435.624236909611
$ws.Range("R10").Value = 21920.6181321865
$ws.Range("S10").Value = 0.2073546078368376
$ws.Range("T10").Value = 0.2073546078368376
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.909391
$ws.Range("H11").Value = 368.728173
$ws.Range("I11").Value = 0.4861499128584522
$ws.Range("J11").Value = 0.4861499128584522
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 18.17573333333333
$ws.Range("N11").Value = 54.5272
$ws.Range("O11").Value = 0.3912102278294916
$ws.Range("P11").Value = 0.3912102278294917
$ws.Range("Q11").Value = 2233.9683149784
$ws.Range("R11").Value = 20105.7148348056
$ws.Range("S11").Value = 0.1901868181686426
$ws.Range("T11").Value = 0.1901868181686426
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 122.909391
$ws.Range("H12").Value = 368.728173
$ws.Range("I12").Value = 0.4861499128584522
$ws.Range("J12").Value = 0.4861499128584522
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.429072666666666
$ws.Range("N12").Value = 19.287218
$ws.Range("O12").Value = 0.1383778545015528
$ws.Range("P12").Value = 0.1383778545015528
$ws.Range("Q12").Value = 790.1934061547458
$ws.Range("R12").Value = 7111.740655392713
$ws.Range("S12").Value = 0.06727238190746948
$ws.Range("T12").Value = 0.06727238190746949
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 122.909391
$ws.Range("H13").Value = 368.728173
$ws.Range("I13").Value = 0.4861499128584522
$ws.Range("J13").Value = 0.4861499128584522
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.039044333333333
$ws.Range("N13").Value = 6.117133
$ws.Range("O13").Value = 0.04388791272233494
$ws.Range("P13").Value = 0.04388791272233494
$ws.Range("Q13").Value = 250.617697232001
$ws.Range("R13").Value = 2255.559275088009
$ws.Range("S13").Value = 0.02133610494550248
$ws.Range("T13").Value = 0.02133610494550249
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.127183
$ws.Range("H14").Value = 45.381549
$ws.Range("I14").Value = 0.05983333443775553
$ws.Range("J14").Value = 0.05983333443775553
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 19.816421
$ws.Range("N14").Value = 59.449263
$ws.Range("O14").Value = 0.4265240049466206
$ws.Range("P14").Value = 0.4265240049466206
$ws.Range("Q14").Value = 299.766626872043
$ws.Range("R14").Value = 2697.899641848387
$ws.Range("S14").Value = 0.02552035343370204
$ws.Range("T14").Value = 0.02552035343370205
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.127183
$ws.Range("H15").Value = 45.381549
$ws.Range("I15").Value = 0.05983333443775553
$ws.Range("J15").Value = 0.05983333443775553
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 18.17573333333333
$ws.Range("N15").Value = 54.5272
$ws.Range("O15").Value = 0.3912102278294916
$ws.Range("P15").Value = 0.3912102278294917
$ws.Range("Q15").Value = 274.9476442925334
$ws.Range("R15").Value = 2474.5287986328
$ws.Range("S15").Value = 0.02340741239719251
$ws.Range("T15").Value = 0.02340741239719251
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.127183
$ws.Range("H16").Value = 45.381549
$ws.Range("I16").Value = 0.05983333443775553
$ws.Range("J16").Value = 0.05983333443775553
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.429072666666666
$ws.Range("N16").Value = 19.287218
$ws.Range("O16").Value = 0.1383778545015528
$ws.Range("P16").Value = 0.1383778545015528
$ws.Range("Q16").Value = 97.25375874896467
$ws.Range("R16").Value = 875.283828740682
$ws.Range("S16").Value = 0.008279608447170484
$ws.Range("T16").Value = 0.008279608447170485
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.127183
$ws.Range("H17").Value = 45.381549
$ws.Range("I17").Value = 0.05983333443775553
$ws.Range("J17").Value = 0.05983333443775553
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.039044333333333
$ws.Range("N17").Value = 6.117133
$ws.Range("O17").Value = 0.04388791272233494
$ws.Range("P17").Value = 0.04388791272233494
$ws.Range("Q17").Value = 30.84499677544634
$ws.Range("R17").Value = 277.604970979017
$ws.Range("S17").Value = 0.002625960159690492
$ws.Range("T17").Value = 0.002625960159690492
